$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 7597
$ws.Range("F4").Value = 7784
$ws.Range("F5").Value = 36
$ws.Range("F8").Value = 6481
$ws.Range("F9").Value = 3332
$ws.Range("F11").Value = 3688
$ws.Range("F12").Value = 39
$ws.Range("F13").Value = 34
$ws.Range("F14").Value = 35
$ws.Range("F15").Value = 53
$ws.Range("F16").Value = 38
$ws.Range("F17").Value = 458
$ws.Range("F19").Value = 305
$ws.Range("F20").Value = 313
$ws.Range("F21").Value = 3768
$ws.Range("F23").Value = 360
$ws.Range("F24").Value = 952
$ws.Range("F25").Value = 277
$ws.Range("F26").Value = 1424
$ws.Range("F27").Value = 75
$ws.Range("F28").Value = 48
$ws.Range("F29").Value = 2719
$ws.Range("F30").Value = 1723
$ws.Range("F32").Value = 38
$ws.Range("F33").Value = 50
$ws.Range("F34").Value = 3556
$ws.Range("F35").Value = 270
$ws.Range("F36").Value = 271
$ws.Range("F38").Value = 915
$ws.Range("F39").Value = 521
$ws.Range("F40").Value = 1378
$ws.Range("F42").Value = 541
$ws.Range("F43").Value = 624

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 407
$ws.Range("F11").Value = 37

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 130

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 130
$ws.Range("F6").Value = 7597
$ws.Range("F7").Value = 7784
$ws.Range("F8").Value = 36
$ws.Range("F11").Value = 6481
$ws.Range("F12").Value = 3332
$ws.Range("F13").Value = 3688
$ws.Range("F14").Value = 34
$ws.Range("F15").Value = 53
$ws.Range("F16").Value = 38
$ws.Range("F17").Value = 458
$ws.Range("F19").Value = 305
$ws.Range("F21").Value = 313
$ws.Range("F22").Value = 3768
$ws.Range("F24").Value = 37
$ws.Range("F26").Value = 360
$ws.Range("F27").Value = 952
$ws.Range("F28").Value = 277
$ws.Range("F29").Value = 1424
$ws.Range("F30").Value = 75
$ws.Range("F31").Value = 48
$ws.Range("F32").Value = 2719
$ws.Range("F33").Value = 1723
$ws.Range("F35").Value = 38
$ws.Range("F36").Value = 50
$ws.Range("F38").Value = 3556
$ws.Range("F39").Value = 270
$ws.Range("F40").Value = 271
$ws.Range("F43").Value = 915
$ws.Range("F44").Value = 521
$ws.Range("F45").Value = 1378
$ws.Range("F48").Value = 541
$ws.Range("F49").Value = 624
